# "latest update from Narjes"
# Fills in the "Type" column (D) on the "RAW Data" sheet for the
# olist_order_payments_dataset, olist_order_reviews_dataset and
# olist_orders_dataset (AMIR) tables, and corrects the data type on the
# order_status row from NVARCHAR(60) to VarChar(15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAW Data")

# olist_order_reviews_dataset - review_comment_title (introduces "VarChar(50)")
$ws.Range("D52").Value = "VarChar(50)"

# olist_order_payments_dataset (rows 40-44)
$ws.Range("D40").Value = "VarChar(40)"   # order_id
$ws.Range("D41").Value = "smallint"      # payment_sequential
$ws.Range("D42").Value = "VarChar(15)"   # payment_type
$ws.Range("D43").Value = "Integer"       # payment_installments
$ws.Range("D44").Value = "real"          # payment_value

# olist_order_reviews_dataset (rows 49-55)
$ws.Range("D49").Value = "VarChar(40)"   # review_id
$ws.Range("D50").Value = "VarChar(40)"   # order_id
$ws.Range("D51").Value = "smallint"      # review_score
$ws.Range("D53").Value = "Text"          # review_comment_message
$ws.Range("D54").Value = "Date"          # review_creation_date
$ws.Range("D55").Value = "Date"          # review_answer_timestamp

# olist_orders_dataset (AMIR) (rows 60-67)
$ws.Range("D60").Value = "VarChar(40)"   # order_id
$ws.Range("D61").Value = "VarChar(40)"   # customer_id
$ws.Range("D62").Value = "VarChar(15)"   # order_status (was NVARCHAR(60))
$ws.Range("D63").Value = "Date"          # order_purchase_timestamp
$ws.Range("D64").Value = "Date"          # order_approved_at
$ws.Range("D65").Value = "Date"          # order_delivered_carrier_date
$ws.Range("D66").Value = "Date"          # order_delivered_customer_date
$ws.Range("D67").Value = "Date"          # order_estimated_delivery_date

# Leave the cursor on the last-edited cell and make "Code" the active tab,
# matching the end-of-session UI state.
$ws.Range("D67").Select()

$wsCode = $wb.Worksheets.Item("Code")
$wsCode.Activate()
$wsCode.Range("M14").Select()
